$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.481.30"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").Value = "1.825.36"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'312.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").Value = "'0.3616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.07198"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("D10").Value = "'0.8609"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.13%  "

$ws.Range("D11").Value = "'20.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("D12").Value = "1.829.90"
$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").Value = "'5.391"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").Value = "'6.483"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").Value = "'0.06929"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "'80.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").Value = "'0.000008858"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").Value = "'15.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").Value = "27.471.07"
$ws.Range("E21").Value = "  -0.91%  "

$ws.Range("D22").Value = "'5.116"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.75%  "

$ws.Range("D23").Value = "'10.93"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "2.052.36"
$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("D25").Value = "'1.986"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").Value = "'155.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "'18.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "

$ws.Range("D28").Value = "'5.138"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.93%  "

$ws.Range("D29").Value = "'114.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.40%  "

$ws.Range("D30").Value = "'1.796"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.91%  "

$ws.Range("D31").Value = "'0.08839"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").Value = "'0.7472"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.86%  "

$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("D34").Value = "'4.533"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("D35").Value = "'1.120"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D37").Value = "'1.086"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.63%  "

$ws.Range("D38").Value = "'0.05280"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.63%  "

$ws.Range("D39").Value = "'0.01917"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").Value = "'2.791"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "

$ws.Range("D41").Value = "'0.5062"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").Value = "'0.1644"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("D43").Value = "'6.436"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.96%  "

$ws.Range("D44").Value = "'8.327"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("D45").Value = "'10.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("D46").Value = "'105.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.47%  "

$ws.Range("D47").Value = "'0.06450"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.42%  "

$ws.Range("D48").Value = "'0.4676"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.71%  "

$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").Value = "'1.612"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.30%  "

$ws.Range("D51").Value = "'63.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
